$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 80 ("Z04_B01_P02_Ib01_I01" / "Schulen mit BNE-Label") needs to move up
# to become row 20, shifting the existing rows 20-79 down by one (to 21-80).
$ws.Rows.Item(80).Cut() | Out-Null
$ws.Rows.Item(20).Insert() | Out-Null
